$wb = $excel.ActiveWorkbook

# --- Create the three new sheets, in final left-to-right order ---
# Worksheets.Add() with no args inserts before the active sheet; passing
# "After" keeps the new sheet right after the one we just added so the
# final left-to-right order matches TestValidLogin, TestInvalidLogin,
# TestValidLoginLogout.
$sValid = $wb.Worksheets.Add()
$sValid.Name = "TestValidLogin"

$sInvalid = $wb.Worksheets.Add($null, $wb.Worksheets.Item("TestValidLogin"))
$sInvalid.Name = "TestInvalidLogin"

$sLogout = $wb.Worksheets.Add($null, $wb.Worksheets.Item("TestInvalidLogin"))
$sLogout.Name = "TestValidLoginLogout"

# The original placeholder sheet is no longer needed - remove it now that
# the workbook has more than one sheet (Excel refuses to delete the last one).
$wb.Worksheets.Item("Sheet1").Delete()

# --- Populate TestValidLogin (admin / manager) ---
# Write row 2 before row 1 so the shared-string table picks up
# admin, manager, Username, Password in that exact order.
$wsValid = $wb.Worksheets.Item("TestValidLogin")
$wsValid.Range("A2").Value = "admin"
$wsValid.Range("B2").Value = "manager"
$wsValid.Range("A1").Value = "Username"
$wsValid.Range("B1").Value = "Password"
$wsValid.Columns("A").ColumnWidth = 10
$wsValid.Columns("B").ColumnWidth = 9.42578125
$wsValid.Range("A1:B2").Select()

# --- Populate TestInvalidLogin (UserName / passowrd, abcd / xyz) ---
$wsInvalid = $wb.Worksheets.Item("TestInvalidLogin")
$wsInvalid.Range("A1").Value = "UserName"
$wsInvalid.Range("B1").Value = "passowrd"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"
$wsInvalid.Range("B3").Select()

# --- Populate TestValidLoginLogout (same data as TestValidLogin) ---
$wsLogout = $wb.Worksheets.Item("TestValidLoginLogout")
$wsLogout.Range("A1").Value = "Username"
$wsLogout.Range("B1").Value = "Password"
$wsLogout.Range("A2").Value = "admin"
$wsLogout.Range("B2").Value = "manager"
$wsLogout.Range("A1:B2").Select()

# Make TestValidLoginLogout the active tab, as in the final workbook.
$wb.Worksheets.Item("TestValidLoginLogout").Activate()
